# Regenerate column G ("K" = strikeouts) values for soria_joakim save_data.
# The workbook stores per-appearance game log data; column G holds the
# newly (re)computed "K" stat (replacing the old Strike# derived figure).
# New values below were produced by the upstream regeneration (std/mean,
# calc and write s_vals) and are applied directly to the sheet cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = [ordered]@{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 2
    13 = 1
    14 = 2
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 2
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 3
    30 = 1
    31 = 1
    32 = 2
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 2
    38 = 1
    39 = 0
    40 = 0
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
